$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table refresh (GitHub Actions scheduled update).
# D column holds price strings that look numeric (e.g. "229.46", "37.345.48")
# but must stay plain text exactly as scraped, so force text format first.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.345.48'
$ws.Range('E2').Value = '  -1.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.047.85'
$ws.Range('E3').Value = '  -1.44%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.46'
$ws.Range('E5').Value = '  -1.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('E6').Value = '  -1.87%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.63'
$ws.Range('E8').Value = '  -3.38%  '

$ws.Range('E9').Value = '  -2.13%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0786'
$ws.Range('E10').Value = '  +0.03%  '

$ws.Range('E11').Value = '  -2.03%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.78'
$ws.Range('E12').Value = '  -0.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.350.18'
$ws.Range('E13').Value = '  -1.41%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.72'
$ws.Range('E14').Value = '  -1.91%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.754'
$ws.Range('E15').Value = '  -3.45%  '

$ws.Range('E16').Value = '  -1.21%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.049.85'
$ws.Range('E17').Value = '  -1.29%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.205.61'
$ws.Range('E18').Value = '  -1.28%  '

$ws.Range('E19').Value = '  -0.97%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.40'
$ws.Range('E20').Value = '  -3.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0827'
$ws.Range('E21').Value = '  -1.89%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.08'
$ws.Range('E22').Value = '  -1.47%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  +0.47%  '

$ws.Range('E25').Value = '  -4.84%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.69'
$ws.Range('E26').Value = '  -0.76%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.23'
$ws.Range('E27').Value = '  -3.22%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.128'
$ws.Range('E28').Value = '  -7.42%  '

$ws.Range('E29').Value = '  -2.11%  '

$ws.Range('E30').Value = '  -3.97%  '

$ws.Range('E31').Value = '  -1.86%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  -4.50%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0616'
$ws.Range('E33').Value = '  -2.83%  '

$ws.Range('E34').Value = '  -1.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.47'
$ws.Range('E35').Value = '  +0.61%  '

$ws.Range('E36').Value = '  +1.29%  '

$ws.Range('E37').Value = '  +0.31%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.23'
$ws.Range('E38').Value = '  -5.26%  '

$ws.Range('E39').Value = '  -3.29%  '

$ws.Range('E40').Value = '  -5.12%  '

$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.484.72'
$ws.Range('E41').Value = '  +2.31%  '

$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.90'
$ws.Range('E42').Value = '  -0.89%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.90'
$ws.Range('E43').Value = '  -0.93%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0944'
$ws.Range('E44').Value = '  -3.40%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.51'
$ws.Range('E45').Value = '  -4.56%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.16'
$ws.Range('E46').Value = '  +0.66%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.97'
$ws.Range('E47').Value = '  -3.16%  '

$ws.Range('E48').Value = '  -4.23%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.10'
$ws.Range('E49').Value = '  -4.03%  '

$ws.Range('E50').Value = '  -2.29%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.235.97'
$ws.Range('E51').Value = '  -1.49%  '

